$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column A to Text format so date-like strings
# ("2025-02-14" etc.) are stored as literal text, not converted to date serials.
$ws.Range("A41:A52").NumberFormat = "@"

$ws.Range("A41").Value = "2025-02-14"
$ws.Range("B41").Value = "sleep"
$ws.Range("C41").Value = $False
$ws.Range("D41").Value = $False
$ws.Range("E41").Value = $True
$ws.Range("F41").Value = $False
$ws.Range("G41").Value = $False
$ws.Range("H41").Value = $True
$ws.Range("I41").Value = $True
$ws.Range("J41").Value = $True
$ws.Range("K41").Value = $False
$ws.Range("L41").Value = $True
$ws.Range("M41").Value = $True
$ws.Range("N41").Value = $True
$ws.Range("O41").Value = $True

$ws.Range("A42").Value = "2025-02-14"
$ws.Range("B42").Value = "activity"
$ws.Range("C42").Value = $True
$ws.Range("D42").Value = $False
$ws.Range("E42").Value = $False
$ws.Range("F42").Value = $True
$ws.Range("G42").Value = $True
$ws.Range("H42").Value = $True
$ws.Range("I42").Value = $False
$ws.Range("J42").Value = $True
$ws.Range("K42").Value = $False
$ws.Range("L42").Value = $True
$ws.Range("M42").Value = $False
$ws.Range("N42").Value = $False
$ws.Range("O42").Value = $False

$ws.Range("A43").Value = "2025-02-14"
$ws.Range("B43").Value = "weekly_activity"
$ws.Range("C43").Value = $True
$ws.Range("D43").Value = $False
$ws.Range("E43").Value = $True
$ws.Range("F43").Value = $True
$ws.Range("G43").Value = $True
$ws.Range("H43").Value = $True
$ws.Range("I43").Value = $True
$ws.Range("J43").Value = $True
$ws.Range("K43").Value = $True
$ws.Range("L43").Value = $True
$ws.Range("M43").Value = $True
$ws.Range("N43").Value = $False
$ws.Range("O43").Value = $False

$ws.Range("A44").Value = "2025-02-15"
$ws.Range("B44").Value = "sleep"
$ws.Range("C44").Value = $True
$ws.Range("D44").Value = $True
$ws.Range("E44").Value = $True
$ws.Range("F44").Value = $True
$ws.Range("G44").Value = $True
$ws.Range("H44").Value = $True
$ws.Range("I44").Value = $False
$ws.Range("J44").Value = $True
$ws.Range("K44").Value = $True
$ws.Range("L44").Value = $True
$ws.Range("M44").Value = $True
$ws.Range("N44").Value = $True
$ws.Range("O44").Value = $True

$ws.Range("A45").Value = "2025-02-15"
$ws.Range("B45").Value = "activity"
$ws.Range("C45").Value = $True
$ws.Range("D45").Value = $False
$ws.Range("E45").Value = $True
$ws.Range("F45").Value = $False
$ws.Range("G45").Value = $True
$ws.Range("H45").Value = $True
$ws.Range("I45").Value = $True
$ws.Range("J45").Value = $True
$ws.Range("K45").Value = $False
$ws.Range("L45").Value = $True
$ws.Range("M45").Value = $False
$ws.Range("N45").Value = $False
$ws.Range("O45").Value = $False

$ws.Range("A46").Value = "2025-02-15"
$ws.Range("B46").Value = "weekly_activity"
$ws.Range("C46").Value = $False
$ws.Range("D46").Value = $False
$ws.Range("E46").Value = $False
$ws.Range("F46").Value = $False
$ws.Range("G46").Value = $False
$ws.Range("H46").Value = $False
$ws.Range("I46").Value = $False
$ws.Range("J46").Value = $False
$ws.Range("K46").Value = $False
$ws.Range("L46").Value = $False
$ws.Range("M46").Value = $False
$ws.Range("N46").Value = $False
$ws.Range("O46").Value = $False

$ws.Range("A47").Value = "2025-02-16"
$ws.Range("B47").Value = "sleep"
$ws.Range("C47").Value = $True
$ws.Range("D47").Value = $False
$ws.Range("E47").Value = $True
$ws.Range("F47").Value = $False
$ws.Range("G47").Value = $False
$ws.Range("H47").Value = $True
$ws.Range("I47").Value = $False
$ws.Range("J47").Value = $True
$ws.Range("K47").Value = $True
$ws.Range("L47").Value = $True
$ws.Range("M47").Value = $True
$ws.Range("N47").Value = $True
$ws.Range("O47").Value = $True

$ws.Range("A48").Value = "2025-02-16"
$ws.Range("B48").Value = "activity"
$ws.Range("C48").Value = $False
$ws.Range("D48").Value = $False
$ws.Range("E48").Value = $False
$ws.Range("F48").Value = $True
$ws.Range("G48").Value = $True
$ws.Range("H48").Value = $True
$ws.Range("I48").Value = $False
$ws.Range("J48").Value = $False
$ws.Range("K48").Value = $False
$ws.Range("L48").Value = $False
$ws.Range("M48").Value = $False
$ws.Range("N48").Value = $False
$ws.Range("O48").Value = $False

$ws.Range("A49").Value = "2025-02-16"
$ws.Range("B49").Value = "weekly_activity"
$ws.Range("C49").Value = $False
$ws.Range("D49").Value = $False
$ws.Range("E49").Value = $False
$ws.Range("F49").Value = $False
$ws.Range("G49").Value = $False
$ws.Range("H49").Value = $False
$ws.Range("I49").Value = $False
$ws.Range("J49").Value = $False
$ws.Range("K49").Value = $False
$ws.Range("L49").Value = $False
$ws.Range("M49").Value = $False
$ws.Range("N49").Value = $False
$ws.Range("O49").Value = $False

$ws.Range("A50").Value = "2025-02-17"
$ws.Range("B50").Value = "sleep"
$ws.Range("C50").Value = $False
$ws.Range("D50").Value = $False
$ws.Range("E50").Value = $True
$ws.Range("F50").Value = $True
$ws.Range("G50").Value = $True
$ws.Range("H50").Value = $True
$ws.Range("I50").Value = $True
$ws.Range("J50").Value = $True
$ws.Range("K50").Value = $False
$ws.Range("L50").Value = $True
$ws.Range("M50").Value = $True
$ws.Range("N50").Value = $True
$ws.Range("O50").Value = $True

$ws.Range("A51").Value = "2025-02-17"
$ws.Range("B51").Value = "activity"
$ws.Range("C51").Value = $False
$ws.Range("D51").Value = $False
$ws.Range("E51").Value = $True
$ws.Range("F51").Value = $True
$ws.Range("G51").Value = $False
$ws.Range("H51").Value = $True
$ws.Range("I51").Value = $True
$ws.Range("J51").Value = $True
$ws.Range("K51").Value = $False
$ws.Range("L51").Value = $True
$ws.Range("M51").Value = $True
$ws.Range("N51").Value = $False
$ws.Range("O51").Value = $False

$ws.Range("A52").Value = "2025-02-17"
$ws.Range("B52").Value = "weekly_activity"
$ws.Range("C52").Value = $False
$ws.Range("D52").Value = $False
$ws.Range("E52").Value = $False
$ws.Range("F52").Value = $False
$ws.Range("G52").Value = $False
$ws.Range("H52").Value = $False
$ws.Range("I52").Value = $False
$ws.Range("J52").Value = $True
$ws.Range("K52").Value = $False
$ws.Range("L52").Value = $False
$ws.Range("M52").Value = $False
$ws.Range("N52").Value = $False
$ws.Range("O52").Value = $False

# Restore default (Normal) style on column A so no stray number-format
# style index lingers on these cells (matches original workbook formatting).
$ws.Range("A41:A52").Style = "Normal"

